# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - row => new value for column F
$sheet1Updates = @{
    2  = 674
    3  = 440
    5  = 1826
    6  = 1435
    8  = 1702
    11 = 640
    12 = 25
    13 = 53
    20 = 60
    22 = 4383
    25 = 92
    26 = 2115
    28 = 1994
}

$wsExpo = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $wsExpo.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# Sheet "全部类型" (fourth sheet) - row => new value for column F
$sheet4Updates = @{
    2  = 674
    3  = 440
    5  = 1826
    6  = 1435
    8  = 1702
    11 = 640
    12 = 25
    13 = 53
    20 = 60
    22 = 4383
    27 = 92
    28 = 2115
    30 = 1994
}

$wsAll = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
